$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row data: row number, word (A), image (B), category (C)
$rows = @(
    @(2, "wenden", "none", "none"),
    @(3, "machen", "house/house027.jpg", "house"),
    @(4, "jubeln", "flower/flower019.jpg", "flower"),
    @(5, "rufen", "none", "none"),
    @(6, "leugnen", "house/house005.jpg", "house"),
    @(7, "retten", "flower/flower012.jpg", "flower"),
    @(8, "enden", "none", "none"),
    @(9, "nerven", "flower/flower001.jpg", "flower"),
    @(10, "fragen", "flower/flower003.jpg", "flower"),
    @(11, "opfern", "none", "none"),
    @(12, "ruhen", "house/house030.jpg", "house"),
    @(13, "achten", "house/house012.jpg", "house"),
    @(14, "schalten", "none", "none"),
    @(15, "stärken", "flower/flower007.jpg", "flower"),
    @(16, "reizen", "house/house031.jpg", "house"),
    @(17, "klagen", "none", "none"),
    @(18, "schwimmen", "house/house019.jpg", "house"),
    @(19, "bitten", "flower/flower022.jpg", "flower"),
    @(20, "sparen", "none", "none"),
    @(21, "runden", "house/house002.jpg", "house"),
    @(22, "öffnen", "flower/flower021.jpg", "flower"),
    @(23, "hören", "none", "none"),
    @(24, "tropfen", "flower/flower030.jpg", "flower"),
    @(25, "treiben", "flower/flower024.jpg", "flower"),
    @(26, "drohen", "none", "none"),
    @(27, "spenden", "house/house008.jpg", "house"),
    @(28, "mühen", "house/house022.jpg", "house"),
    @(29, "orten", "none", "none"),
    @(30, "lügen", "house/house016.jpg", "house"),
    @(31, "gleichen", "house/house014.jpg", "house"),
    @(32, "weigern", "none", "none"),
    @(33, "stopfen", "flower/flower004.jpg", "flower"),
    @(34, "bellen", "house/house010.jpg", "house"),
    @(35, "ächzen", "none", "none"),
    @(36, "wüten", "flower/flower014.jpg", "flower"),
    @(37, "altern", "flower/flower010.jpg", "flower"),
    @(38, "kosten", "none", "none"),
    @(39, "mauern", "flower/flower026.jpg", "flower"),
    @(40, "schleppen", "house/house028.jpg", "house"),
    @(41, "dauern", "none", "none"),
    @(42, "testen", "house/house009.jpg", "house"),
    @(43, "wagen", "flower/flower025.jpg", "flower"),
    @(44, "stören", "none", "none"),
    @(45, "kommen", "flower/flower015.jpg", "flower"),
    @(46, "posten", "house/house001.jpg", "house"),
    @(47, "bremsen", "none", "none"),
    @(48, "leiten", "flower/flower033.jpg", "flower"),
    @(49, "segeln", "house/house006.jpg", "house")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
